$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.690.66'
$ws.Range("E2").Value = '  -5.74%  '
$ws.Range("D3").Value = '3.319.89'
$ws.Range("E3").Value = '  -6.50%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.06%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.69%  '
$ws.Range("D9").Value = '3.314.93'
$ws.Range("E9").Value = '  -6.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.186'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.586'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.37'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000264'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -8.66%  '
$ws.Range("D14").Value = '3.862.15'
$ws.Range("E14").Value = '  -6.16%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '602.81'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -9.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.63%  '
$ws.Range("D18").Value = '65.681.83'
$ws.Range("E18").Value = '  -5.85%  '
$ws.Range("E19").Value = '  -3.97%  '
$ws.Range("D20").Value = '3.318.45'
$ws.Range("E20").Value = '  -6.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -9.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.902'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '102.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.33%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -9.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("E28").Value = '  -9.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.48'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.44%  '
$ws.Range("E32").Value = '  -11.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.82%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '3.827.34'
$ws.Range("E35").Value = '  +2.31%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.104'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '523.85'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.45%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '56.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.15%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.41%  '
$ws.Range("D41").Value = '0.0₃0708'
$ws.Range("E41").Value = '  -13.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.125'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.96%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.338'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.14%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '31.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.97%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.96%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0410'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.58%  '
$ws.Range("B48").Value = 'CoreDAO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +14.10%  '
$ws.Range("E49").Value = '  -5.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.56%  '
$ws.Range("E51").Value = '  -0.04%  '
